$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new work-log entry on row 13 (Laurent's Matlab/video-processing task)
$ws.Range("A13").Value = "Laurent"
$ws.Range("B13").Value = "Introduction to video processing with Matlab"
$ws.Range("C13").Value = "Learning to use Matlab for video processing"
$ws.Range("D13").Value = 42450
$ws.Range("D13").NumberFormat = "d-mmm"
$ws.Range("E13").Value = "1h"

# Move the active selection to A14, matching where the user ended up after filling the row
$ws.Range("A14").Select() | Out-Null
